$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value2 = 13009.223
$ws.Range("I18").Value2 = 15428.714
$ws.Range("K18").Value2 = 15428.714
$ws.Range("M18").Value2 = -15144.714

$ws.Range("H33").Value2 = 258.33334
$ws.Range("I33").Value2 = 287.5
$ws.Range("J33").Value2 = 200
$ws.Range("K33").Value2 = 287.5
$ws.Range("L33").Value2 = 200
$ws.Range("M33").Value2 = -58.5
$ws.Range("N33").Value2 = -658

$ws.Range("H41").Value2 = 679.3889
$ws.Range("I41").Value2 = 842
$ws.Range("J41").Value2 = 354.16666
$ws.Range("K41").Value2 = 842
$ws.Range("L41").Value2 = 354.16666
$ws.Range("M41").Value2 = -402
$ws.Range("N41").Value2 = -1234.16666

$ws.Range("H53").Value2 = 126.4
$ws.Range("I53").Value2 = 102.5
$ws.Range("J53").Value2 = 142.33333
$ws.Range("K53").Value2 = 102.5
$ws.Range("L53").Value2 = 142.33333
$ws.Range("M53").Value2 = 534.5
$ws.Range("N53").Value2 = -1416.33333

$ws.Range("H76").Value2 = 3421.1177
$ws.Range("I76").Value2 = 3432.9285
$ws.Range("J76").Value2 = 3366
$ws.Range("K76").Value2 = 3432.9285
$ws.Range("L76").Value2 = 3366
$ws.Range("M76").Value2 = -3117.9285
$ws.Range("N76").Value2 = -3996

$ws.Range("H79").Value2 = 3421.1177
$ws.Range("I79").Value2 = 3432.9285
$ws.Range("J79").Value2 = 3366
$ws.Range("K79").Value2 = 3432.9285
$ws.Range("L79").Value2 = 3366
$ws.Range("M79").Value2 = -2340.9285
$ws.Range("N79").Value2 = -5550

$ws.Range("H81").Value2 = 0
$ws.Range("J81").Value2 = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("N84").ClearContents()

$ws.Range("H86").Value2 = 1314.4
$ws.Range("I86").Value2 = 893
$ws.Range("K86").Value2 = 893
$ws.Range("M86").Value2 = 230

$ws.Range("H89").Value2 = 1314.4
$ws.Range("I89").Value2 = 893
$ws.Range("K89").Value2 = 4465
$ws.Range("M89").Value2 = 1151

$ws.Range("H92").Value2 = 573.2857
$ws.Range("I92").Value2 = 500.75
$ws.Range("K92").Value2 = 500.75
$ws.Range("M92").Value2 = 747.25

$ws.Range("H98").Value2 = 2263.25
$ws.Range("I98").Value2 = 1184.5834
$ws.Range("J98").Value2 = 5499.25
$ws.Range("K98").Value2 = 1184.5834
$ws.Range("L98").Value2 = 5499.25
$ws.Range("M98").Value2 = 313.4166
$ws.Range("N98").Value2 = -8495.25

$ws.Range("H105").Value2 = 40671
$ws.Range("J105").Value2 = 40671
$ws.Range("L105").Value2 = 40671
$ws.Range("N105").Value2 = -47659

$ws.Range("H122").Value2 = 2263.25
$ws.Range("I122").Value2 = 1184.5834
$ws.Range("J122").Value2 = 5499.25
$ws.Range("K122").Value2 = 3553.7502
$ws.Range("L122").Value2 = 16497.75
$ws.Range("M122").Value2 = -1103.7502
$ws.Range("N122").Value2 = -21397.75

$ws.Range("H132").Value2 = 3820.0334
$ws.Range("I132").Value2 = 2127.423
$ws.Range("J132").Value2 = 14822
$ws.Range("K132").Value2 = 6382.268999999999
$ws.Range("L132").Value2 = 44466
$ws.Range("M132").Value2 = -3852.268999999999
$ws.Range("N132").Value2 = -49526

$ws.Range("H138").Value2 = 2308.647
$ws.Range("J138").Value2 = 2623.4443
$ws.Range("L138").Value2 = 7870.3329
$ws.Range("N138").Value2 = -18150.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value2 = 777.3333
$ws.Range("I6").Value2 = 777.3333
$ws.Range("K6").Value2 = 777.3333
$ws.Range("M6").Value2 = -604.3333

$ws.Range("H32").Value2 = 2201.4583
$ws.Range("I32").Value2 = 1971.1086
$ws.Range("J32").Value2 = 7499.5
$ws.Range("K32").Value2 = 1971.1086
$ws.Range("L32").Value2 = 7499.5
$ws.Range("M32").Value2 = -1684.1086
$ws.Range("N32").Value2 = -8073.5

$ws.Range("H61").Value2 = 2082.5
$ws.Range("I61").Value2 = 2039
$ws.Range("J61").Value2 = 2126
$ws.Range("K61").Value2 = 2039
$ws.Range("L61").Value2 = 2126
$ws.Range("M61").Value2 = -1827
$ws.Range("N61").Value2 = -2550

$ws.Range("H110").Value2 = 965.0526
$ws.Range("I110").Value2 = 965.0526
$ws.Range("K110").Value2 = 965.0526
$ws.Range("M110").Value2 = 1079.9474

$ws.Range("H130").Value2 = 37500
$ws.Range("J130").Value2 = 37500
$ws.Range("L130").Value2 = 37500
$ws.Range("N130").Value2 = -47540

$ws.Range("H132").Value2 = 1658.3715
$ws.Range("I132").Value2 = 1438.1666
$ws.Range("K132").Value2 = 4314.4998
$ws.Range("M132").Value2 = -1784.4998

$ws.Range("H136").Value2 = 2082.5
$ws.Range("I136").Value2 = 2039
$ws.Range("J136").Value2 = 2126
$ws.Range("K136").Value2 = 6117
$ws.Range("L136").Value2 = 6378
$ws.Range("M136").Value2 = -3567
$ws.Range("N136").Value2 = -11478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 1473.96
$ws.Range("I134").Value2 = 1297.7142
$ws.Range("J134").Value2 = 2399.25
$ws.Range("K134").Value2 = 3893.1426
$ws.Range("L134").Value2 = 7197.75
$ws.Range("M134").Value2 = -1358.1426
$ws.Range("N134").Value2 = -12267.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 9802.5
$ws.Range("I31").Value2 = 2877.389
$ws.Range("K31").Value2 = 2877.389
$ws.Range("M31").Value2 = -2582.389

$ws.Range("H34").Value2 = 9802.5
$ws.Range("I34").Value2 = 2877.389
$ws.Range("K34").Value2 = 2877.389
$ws.Range("M34").Value2 = -2675.389

$ws.Range("H37").Value2 = 0
$ws.Range("I37").Value2 = 0
$ws.Range("K37").Value2 = 0
$ws.Range("M37").ClearContents()

$ws.Range("H38").Value2 = 14148.286
$ws.Range("I38").Value2 = 9012.666999999999
$ws.Range("J38").Value2 = 18000
$ws.Range("K38").Value2 = 9012.666999999999
$ws.Range("L38").Value2 = 18000
$ws.Range("M38").Value2 = -8635.666999999999
$ws.Range("N38").Value2 = -18754

$ws.Range("H46").Value2 = 14148.286
$ws.Range("I46").Value2 = 9012.666999999999
$ws.Range("J46").Value2 = 18000
$ws.Range("K46").Value2 = 9012.666999999999
$ws.Range("L46").Value2 = 18000
$ws.Range("M46").Value2 = -8801.666999999999
$ws.Range("N46").Value2 = -18422

$ws.Range("H109").Value2 = 15018.174
$ws.Range("J109").Value2 = 15018.174
$ws.Range("L109").Value2 = 15018.174
$ws.Range("N109").Value2 = -17098.174

$ws.Range("H132").Value2 = 5223.2144
$ws.Range("I132").Value2 = 5298.8184
$ws.Range("J132").Value2 = 4946
$ws.Range("K132").Value2 = 15896.4552
$ws.Range("L132").Value2 = 14838
$ws.Range("M132").Value2 = -13366.4552
$ws.Range("N132").Value2 = -19898

$ws.Range("H133").Value2 = 113666.336
$ws.Range("J133").Value2 = 113666.336
$ws.Range("L133").Value2 = 113666.336
$ws.Range("N133").Value2 = -118726.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 373.33334
$ws.Range("I2").Value2 = 980.6667
$ws.Range("J2").Value2 = 69.666664
$ws.Range("K2").Value2 = 5884.0002
$ws.Range("L2").Value2 = 417.999984
$ws.Range("M2").Value2 = -5771.0002
$ws.Range("N2").Value2 = -643.999984

$ws.Range("H3").Value2 = 3961
$ws.Range("J3").Value2 = 4516.5
$ws.Range("L3").Value2 = 13549.5
$ws.Range("N3").Value2 = -13773.5

$ws.Range("H50").Value2 = 1875474.6
$ws.Range("I50").Value2 = 472.25
$ws.Range("J50").Value2 = 3750477
$ws.Range("K50").Value2 = 1416.75
$ws.Range("L50").Value2 = 11251431
$ws.Range("M50").Value2 = -935.75
$ws.Range("N50").Value2 = -11252393

$ws.Range("H53").Value2 = 1875474.6
$ws.Range("I53").Value2 = 472.25
$ws.Range("J53").Value2 = 3750477
$ws.Range("K53").Value2 = 1416.75
$ws.Range("L53").Value2 = 11251431
$ws.Range("M53").Value2 = -935.75
$ws.Range("N53").Value2 = -11252393

$ws.Range("H102").Value2 = 12599.8
$ws.Range("J102").Value2 = 13999.75
$ws.Range("L102").Value2 = 41999.25
$ws.Range("N102").Value2 = -46867.25

$ws.Range("H131").Value2 = 1906.6154
$ws.Range("I131").Value2 = 3015
$ws.Range("K131").Value2 = 9045
$ws.Range("M131").Value2 = -4005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 6567.8438
$ws.Range("I70").Value2 = 6505
$ws.Range("K70").Value2 = 6505
$ws.Range("M70").Value2 = -6235

$ws.Range("H73").Value2 = 6567.8438
$ws.Range("I73").Value2 = 6505
$ws.Range("K73").Value2 = 6505
$ws.Range("M73").Value2 = -5569

$ws.Range("H102").Value2 = 3029.682
$ws.Range("I102").Value2 = 3032.9
$ws.Range("K102").Value2 = 3032.9
$ws.Range("M102").Value2 = -1410.9

$ws.Range("H132").Value2 = 8468.444
$ws.Range("I132").Value2 = 9883.666999999999
$ws.Range("J132").Value2 = 5638
$ws.Range("K132").Value2 = 29651.001
$ws.Range("L132").Value2 = 16914
$ws.Range("M132").Value2 = -27121.001
$ws.Range("N132").Value2 = -21974

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value2 = 4953.375
$ws.Range("I136").Value2 = 3664
$ws.Range("J136").Value2 = 5727
$ws.Range("K136").Value2 = 10992
$ws.Range("L136").Value2 = 17181
$ws.Range("M136").Value2 = -8442
$ws.Range("N136").Value2 = -22281

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value2 = 25130
$ws.Range("J75").Value2 = 25130
$ws.Range("L75").Value2 = 25130
$ws.Range("N75").Value2 = -27002

$ws.Range("H78").Value2 = 25130
$ws.Range("J78").Value2 = 25130
$ws.Range("L78").Value2 = 75390
$ws.Range("N78").Value2 = -84750

$ws.Range("H109").Value2 = 35000
$ws.Range("J109").Value2 = 35000
$ws.Range("L109").Value2 = 35000
$ws.Range("N109").Value2 = -37774

$ws.Range("H122").Value2 = 1216.8
$ws.Range("I122").Value2 = 1175.579
$ws.Range("J122").Value2 = 2000
$ws.Range("K122").Value2 = 3526.737
$ws.Range("L122").Value2 = 6000
$ws.Range("M122").Value2 = -1076.737
$ws.Range("N122").Value2 = -10900

$ws.Range("H136").Value2 = 3388.4849
$ws.Range("J136").Value2 = 4000
$ws.Range("L136").Value2 = 12000
$ws.Range("N136").Value2 = -17100
